$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("D2").Value2 = "FAPs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 1.553094
$ws.Range("H2").Value2 = 4.659282
$ws.Range("I2").Value2 = 0.6859765954652609
$ws.Range("J2").Value2 = 0.6859765954652609
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 1.036595333333333
$ws.Range("N2").Value2 = 3.109786
$ws.Range("O2").Value2 = 0.393072250513715
$ws.Range("P2").Value2 = 0.393072250513715
$ws.Range("Q2").Value2 = 1.609929992628
$ws.Range("R2").Value2 = 14.489369933652
$ws.Range("S2").Value2 = 0.2696383641792664
$ws.Range("T2").Value2 = 0.2696383641792664

# --- Row 3 updates ---
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 1.553094
$ws.Range("H3").Value2 = 4.659282
$ws.Range("I3").Value2 = 0.6859765954652609
$ws.Range("J3").Value2 = 0.6859765954652609
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 1.600567
$ws.Range("N3").Value2 = 4.801701
$ws.Range("O3").Value2 = 0.6069277494862849
$ws.Range("P3").Value2 = 0.6069277494862849
$ws.Range("Q3").Value2 = 2.485831004298
$ws.Range("R3").Value2 = 22.372479038682
$ws.Range("S3").Value2 = 0.4163382312859945
$ws.Range("T3").Value2 = 0.4163382312859945

# --- Row 4 (new) ---
$ws.Range("A4").Value2 = "sCs"
$ws.Range("B4").Value2 = "Gm13306"
$ws.Range("C4").Value2 = "Ccr10"
$ws.Range("D4").Value2 = "FAPs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 0.7109686666666667
$ws.Range("H4").Value2 = 2.132906
$ws.Range("I4").Value2 = 0.314023404534739
$ws.Range("J4").Value2 = 0.314023404534739
$ws.Range("K4").Value2 = 2
$ws.Range("L4").Value2 = 0.6666666666666666
$ws.Range("M4").Value2 = 1.036595333333333
$ws.Range("N4").Value2 = 3.109786
$ws.Range("O4").Value2 = 0.393072250513715
$ws.Range("P4").Value2 = 0.393072250513715
$ws.Range("Q4").Value2 = 0.7369868020128889
$ws.Range("R4").Value2 = 6.632881218116001
$ws.Range("S4").Value2 = 0.1234338863344486
$ws.Range("T4").Value2 = 0.1234338863344486

# --- Row 5 (new) ---
$ws.Range("A5").Value2 = "sCs"
$ws.Range("B5").Value2 = "Gm13306"
$ws.Range("C5").Value2 = "Ccr10"
$ws.Range("D5").Value2 = "sCs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 0.7109686666666667
$ws.Range("H5").Value2 = 2.132906
$ws.Range("I5").Value2 = 0.314023404534739
$ws.Range("J5").Value2 = 0.314023404534739
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 1.600567
$ws.Range("N5").Value2 = 4.801701
$ws.Range("O5").Value2 = 0.6069277494862849
$ws.Range("P5").Value2 = 0.6069277494862849
$ws.Range("Q5").Value2 = 1.137952985900667
$ws.Range("R5").Value2 = 10.241576873106
$ws.Range("S5").Value2 = 0.1905895182002904
$ws.Range("T5").Value2 = 0.1905895182002904
